# The commit swaps the OOXML content of ppt/theme/theme1.xml (previously the
# "Integral" theme used by the slide master / all slides) and
# ppt/theme/theme2.xml (previously the default "Office Theme" used by the
# notes master), so theme1.xml ends up holding the "Office Theme" palette.
#
# The font scheme (fontScheme) and format scheme (fmtScheme: fills, lines,
# effects, background fills) are byte-for-byte identical between the two
# themes, so the only observable difference is the 12 colour-scheme (clrScheme)
# entries. Re-point every theme colour on the presentation's slide master to
# the "Office Theme" palette to reproduce that swap.
# (RGB values below are plain OLE COLOR ints: R + G*256 + B*65536.)

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0         # dk1      -> 000000
$cs.Item(2).RGB  = 16777215  # lt1      -> FFFFFF
$cs.Item(3).RGB  = 6968388   # dk2      -> 44546A
$cs.Item(4).RGB  = 15132391  # lt2      -> E7E6E6
$cs.Item(5).RGB  = 13998939  # accent1  -> 5B9BD5
$cs.Item(6).RGB  = 3243501   # accent2  -> ED7D31
$cs.Item(7).RGB  = 10855845  # accent3  -> A5A5A5
$cs.Item(8).RGB  = 49407     # accent4  -> FFC000
$cs.Item(9).RGB  = 12874308  # accent5  -> 4472C4
$cs.Item(10).RGB = 4697456   # accent6  -> 70AD47
$cs.Item(11).RGB = 12673797  # hlink    -> 0563C1
$cs.Item(12).RGB = 7491477   # folHlink -> 954F72
